$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82 (shifts existing rows 82..190 down to 83..191)
$ws.Rows("82:82").Insert()

# Populate the newly inserted row 82 with the new weekly price record
$ws.Cells.Item(82, 1).Value = 3
$ws.Cells.Item(82, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(82, 3).Value = "Coquimbo"
$ws.Cells.Item(82, 4).Value = 44467
$ws.Cells.Item(82, 5).Value = 5
$ws.Cells.Item(82, 6).Value = 100112012
$ws.Cells.Item(82, 7).Value = "Espinaca"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 160
$ws.Cells.Item(82, 11).Value = 2500
$ws.Cells.Item(82, 12).Value = 2500
$ws.Cells.Item(82, 13).Value = 2500
$ws.Cells.Item(82, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(82, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(82, 16).Value = 833
$ws.Cells.Item(82, 17).Value = 3
$ws.Cells.Item(82, 18).Value = "Hortaliza"
